$d = $word.ActiveDocument
$wdParagraph = 4

function Get-ParaIndexByText($text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    return $rng.Paragraphs.Item(1).Index
}

function Set-ParaXml($text, $xml) {
    $idx = Get-ParaIndexByText($text)
    $para = $d.Paragraphs.Item($idx)
    $null = $para.Range.InsertXML($xml)
    return $idx
}

# 1. Heading: merge "Agendapunten" / " 3-5-2016" runs (drop proofErr marks) into one run.
$null = Set-ParaXml "Agendapunten" "<w:p><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Agendapunten 3-5-2016</w:t></w:r></w:p>"

# 2. "Writing design report": merge the four split runs into one (same sz/szCs formatting).
$null = Set-ParaXml "riting d" "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t>Writing design report</w:t></w:r></w:p>"

# 3. "Testing basic functionality robots.": merge only the "unctionality" + " " runs.
$null = Set-ParaXml "Testing b" @"
<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t>Testing b</w:t></w:r><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t xml:space='preserve'>asic </w:t></w:r><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t>f</w:t></w:r><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t xml:space='preserve'>unctionality </w:t></w:r><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t xml:space='preserve'>obots. </w:t></w:r></w:p>
"@

# 4. "Developing main algorithm": merge "D" + "eveloping main algorithm" runs.
$null = Set-ParaXml "eveloping main algorithm" "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:rPr><w:sz w:val='21'/><w:szCs w:val='21'/></w:rPr><w:t>Developing main algorithm</w:t></w:r></w:p>"

# 5. "Distribution of tasks" becomes "Peer review design report" at ilvl 1, carrying the
#    _GoBack bookmark, followed by a brand-new "Distribution of tasks" bullet at ilvl 0.
$null = Set-ParaXml "Distribution of tasks" @"
<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Peer review design report</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p><w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Distribution of tasks</w:t></w:r></w:p>
"@

# 6. "AOB" keeps its own bullet but loses the bookmark (it now lives on "Peer review design report").
$null = Set-ParaXml "AOB" "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>AOB</w:t></w:r></w:p>"
